$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) updates to column F ("想去人数" / interest count)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 597
$ws1.Range("F9").Value = 199
$ws1.Range("F11").Value = 476
$ws1.Range("F12").Value = 1430
$ws1.Range("F14").Value = 129
$ws1.Range("F15").Value = 290
$ws1.Range("F17").Value = 102
$ws1.Range("F18").Value = 682
$ws1.Range("F19").Value = 1028
$ws1.Range("F20").Value = 44
$ws1.Range("F21").Value = 271
$ws1.Range("F22").Value = 37
$ws1.Range("F23").Value = 6147
$ws1.Range("F25").Value = 130
$ws1.Range("F26").Value = 127
$ws1.Range("F28").Value = 14899
$ws1.Range("F29").Value = 1478
$ws1.Range("F32").Value = 94
$ws1.Range("F33").Value = 10864
$ws1.Range("F34").Value = 685
$ws1.Range("F35").Value = 4248
$ws1.Range("F36").Value = 188
$ws1.Range("F38").Value = 117

# Sheet "全部类型" (sheet4) updates mirroring the same rows (aggregated sheet)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 597
$ws4.Range("F9").Value = 199
$ws4.Range("F11").Value = 476
$ws4.Range("F12").Value = 1430
$ws4.Range("F14").Value = 129
$ws4.Range("F15").Value = 290
$ws4.Range("F18").Value = 102
$ws4.Range("F19").Value = 682
$ws4.Range("F21").Value = 1028
$ws4.Range("F22").Value = 44
$ws4.Range("F23").Value = 271
$ws4.Range("F24").Value = 37
$ws4.Range("F26").Value = 6147
$ws4.Range("F28").Value = 130
$ws4.Range("F29").Value = 127
$ws4.Range("F31").Value = 14899
$ws4.Range("F32").Value = 1478
$ws4.Range("F35").Value = 94
$ws4.Range("F36").Value = 10864
$ws4.Range("F37").Value = 685
$ws4.Range("F38").Value = 4248
$ws4.Range("F39").Value = 188
$ws4.Range("F41").Value = 117
